$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the direct (border/fill/numberformat) cell formatting from row 12 down to the
# new row 13, so the new cells end up on the same style records (s="3"/s="4"),
# matching how every other data row in the sheet is formatted.
$ws.Range("A12:B12").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122)  # xlPasteFormats

# A13 holds a day-month-year looking string ("04-10-2025"). Assigning it straight to
# .Value would make Excel "smart" parse it into a date serial. Build it as a text
# formula result on a scratch cell first (guarantees text type), then paste only the
# *value* into A13 - this keeps it a shared string, same as the other date cells,
# without touching A13's number format / style.
$ws.Range("Z1").Formula = "=""04-10-2025"""
$ws.Range("Z1").Copy()
$ws.Range("A13").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").Clear()

$ws.Range("B13").Value = "The price of gold in India today is ₹11,940 per gram for 24 karat gold, ₹10,945 per gram for 22 karat gold and ₹8,955 per gram for 18 karat gold (also called 999 gold)."
